# Scheduled market-data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) for the rows whose market snapshot changed on this run.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H28").Value = 1096.1333
$ws_ALC.Range("I28").Value = 1086.4445
$ws_ALC.Range("J28").Value = 1110.6666
$ws_ALC.Range("K28").Value = 1086.4445
$ws_ALC.Range("L28").Value = 1110.6666
$ws_ALC.Range("M28").Value = -601.4445000000001
$ws_ALC.Range("N28").Value = -2080.6666

$ws_ALC.Range("H80").Value = 4084.8484
$ws_ALC.Range("I80").Value = 651.7143
$ws_ALC.Range("J80").Value = 6614.5264
$ws_ALC.Range("K80").Value = 1955.1429
$ws_ALC.Range("L80").Value = 19843.5792
$ws_ALC.Range("M80").Value = -957.1428999999998
$ws_ALC.Range("N80").Value = -21839.5792

$ws_ALC.Range("H83").Value = 4084.8484
$ws_ALC.Range("I83").Value = 651.7143
$ws_ALC.Range("J83").Value = 6614.5264
$ws_ALC.Range("K83").Value = 5865.428699999999
$ws_ALC.Range("L83").Value = 59530.7376
$ws_ALC.Range("M83").Value = -873.4286999999995
$ws_ALC.Range("N83").Value = -69514.73759999999

$ws_ALC.Range("H100").Value = 1953.1428
$ws_ALC.Range("I100").Value = 1953.1428
$ws_ALC.Range("K100").Value = 1953.1428
$ws_ALC.Range("M100").Value = -1412.1428

$ws_ALC.Range("H125").Value = 21752.066
$ws_ALC.Range("J125").Value = 1488.1
$ws_ALC.Range("L125").Value = 13392.9
$ws_ALC.Range("N125").Value = -18312.9

$ws_ALC.Range("H137").Value = 21599.2
$ws_ALC.Range("I137").Value = 32548.2
$ws_ALC.Range("J137").Value = 17219.6
$ws_ALC.Range("K137").Value = 97644.60000000001
$ws_ALC.Range("L137").Value = 51658.8
$ws_ALC.Range("M137").Value = -95094.60000000001
$ws_ALC.Range("N137").Value = -56758.8

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H5").Value = 10829.6
$ws_ARM.Range("I5").Value = 10829.6
$ws_ARM.Range("K5").Value = 10829.6
$ws_ARM.Range("M5").Value = -10717.6

$ws_ARM.Range("H34").Value = 133335.67
$ws_ARM.Range("J34").Value = 100002.55
$ws_ARM.Range("L34").Value = 100002.55
$ws_ARM.Range("N34").Value = -100544.55

$ws_ARM.Range("H45").Value = 5986.143
$ws_ARM.Range("I45").Value = 5400.6
$ws_ARM.Range("K45").Value = 5400.6
$ws_ARM.Range("M45").Value = -5023.6

$ws_ARM.Range("H132").Value = 2881.3
$ws_ARM.Range("I132").Value = 2275.2
$ws_ARM.Range("K132").Value = 6825.599999999999
$ws_ARM.Range("M132").Value = -4295.599999999999

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H4").Value = 10829.6
$ws_BSM.Range("I4").Value = 10829.6
$ws_BSM.Range("K4").Value = 10829.6
$ws_BSM.Range("M4").Value = -10714.6

$ws_BSM.Range("H20").Value = 9502.666999999999
$ws_BSM.Range("I20").Value = 11591.105
$ws_BSM.Range("K20").Value = 11591.105
$ws_BSM.Range("M20").Value = -11344.105

$ws_BSM.Range("H22").Value = 761.5833
$ws_BSM.Range("I22").Value = 758.0909
$ws_BSM.Range("K22").Value = 758.0909
$ws_BSM.Range("M22").Value = -585.0909

$ws_BSM.Range("H105").Value = 1885.0667
$ws_BSM.Range("I105").Value = 1582.6296
$ws_BSM.Range("J105").Value = 4607
$ws_BSM.Range("K105").Value = 1582.6296
$ws_BSM.Range("L105").Value = 4607
$ws_BSM.Range("M105").Value = 164.3704
$ws_BSM.Range("N105").Value = -8101

$ws_BSM.Range("H134").Value = 2574.2163
$ws_BSM.Range("I134").Value = 2335.394
$ws_BSM.Range("K134").Value = 7006.181999999999
$ws_BSM.Range("M134").Value = -4471.181999999999

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 7149694.5
$ws_CRP.Range("I31").Value = 33344300
$ws_CRP.Range("J31").Value = 5710.8184
$ws_CRP.Range("K31").Value = 33344300
$ws_CRP.Range("L31").Value = 5710.8184
$ws_CRP.Range("M31").Value = -33344005
$ws_CRP.Range("N31").Value = -6300.8184

$ws_CRP.Range("H34").Value = 7149694.5
$ws_CRP.Range("I34").Value = 33344300
$ws_CRP.Range("J34").Value = 5710.8184
$ws_CRP.Range("K34").Value = 33344300
$ws_CRP.Range("L34").Value = 5710.8184
$ws_CRP.Range("M34").Value = -33344098
$ws_CRP.Range("N34").Value = -6114.8184

$ws_CRP.Range("H125").Value = 95329.664
$ws_CRP.Range("J125").Value = 95329.664
$ws_CRP.Range("L125").Value = 95329.664
$ws_CRP.Range("N125").Value = -100249.664

$ws_CRP.Range("H132").Value = 84923.586
$ws_CRP.Range("I132").Value = 112407.664
$ws_CRP.Range("J132").Value = 2471.3333
$ws_CRP.Range("K132").Value = 337222.992
$ws_CRP.Range("L132").Value = 7413.999899999999
$ws_CRP.Range("M132").Value = -334692.992
$ws_CRP.Range("N132").Value = -12473.9999

$ws_CRP.Range("H134").Value = 3730.0833
$ws_CRP.Range("I134").Value = 3473
$ws_CRP.Range("J134").Value = 4244.25
$ws_CRP.Range("K134").Value = 10419
$ws_CRP.Range("L134").Value = 12732.75
$ws_CRP.Range("M134").Value = -7884
$ws_CRP.Range("N134").Value = -17802.75

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H36").Value = 11000
$ws_GSM.Range("J36").Value = 11000
$ws_GSM.Range("L36").Value = 11000
$ws_GSM.Range("N36").Value = -11970

$ws_GSM.Range("H132").Value = 4271
$ws_GSM.Range("I132").Value = 4006.5
$ws_GSM.Range("K132").Value = 12019.5
$ws_GSM.Range("M132").Value = -9489.5

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 4226.636
$ws_LTW.Range("I7").Value = 3610.3333
$ws_LTW.Range("K7").Value = 3610.3333
$ws_LTW.Range("M7").Value = -3498.3333

$ws_LTW.Range("H40").Value = 4475.4
$ws_LTW.Range("I40").Value = 4336.2856
$ws_LTW.Range("J40").Value = 4800
$ws_LTW.Range("K40").Value = 4336.2856
$ws_LTW.Range("L40").Value = 4800
$ws_LTW.Range("M40").Value = -4200.2856
$ws_LTW.Range("N40").Value = -5072

$ws_LTW.Range("H46").Value = 4961.222
$ws_LTW.Range("I46").Value = 945.5
$ws_LTW.Range("J46").Value = 8173.8
$ws_LTW.Range("K46").Value = 945.5
$ws_LTW.Range("L46").Value = 8173.8
$ws_LTW.Range("M46").Value = -757.5
$ws_LTW.Range("N46").Value = -8549.799999999999

$ws_LTW.Range("H55").Value = 3206.818
$ws_LTW.Range("I55").Value = 1979.3334
$ws_LTW.Range("K55").Value = 1979.3334
$ws_LTW.Range("M55").Value = -1806.3334

$ws_LTW.Range("H122").Value = 3650.889
$ws_LTW.Range("I122").Value = 3550.4285
$ws_LTW.Range("J122").Value = 4002.5
$ws_LTW.Range("K122").Value = 10651.2855
$ws_LTW.Range("L122").Value = 12007.5
$ws_LTW.Range("M122").Value = -8201.2855
$ws_LTW.Range("N122").Value = -16907.5

$ws_LTW.Range("H126").Value = 4226.636
$ws_LTW.Range("I126").Value = 3610.3333
$ws_LTW.Range("K126").Value = 10830.9999
$ws_LTW.Range("M126").Value = -8360.999899999999

$ws_LTW.Range("H132").Value = 3943.5715
$ws_LTW.Range("I132").Value = 3830.9092
$ws_LTW.Range("K132").Value = 11492.7276
$ws_LTW.Range("M132").Value = -8962.7276

$ws_LTW.Range("H136").Value = 6211.3335
$ws_LTW.Range("I136").Value = 6549.6665
$ws_LTW.Range("K136").Value = 19648.9995
$ws_LTW.Range("M136").Value = -17098.9995

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H15").Value = 29969
$ws_WVR.Range("J15").Value = 29969
$ws_WVR.Range("L15").Value = 29969
$ws_WVR.Range("N15").Value = -30545

$ws_WVR.Range("H40").Value = 60000
$ws_WVR.Range("I40").Value = 60000
$ws_WVR.Range("K40").Value = 60000
$ws_WVR.Range("M40").Value = -59851

$ws_WVR.Range("H54").Value = 4900
$ws_WVR.Range("I54").Value = 4900
$ws_WVR.Range("K54").Value = 4900
$ws_WVR.Range("M54").Value = -4380

$ws_WVR.Range("H64").Value = 0
$ws_WVR.Range("J64").Value = 0
$ws_WVR.Range("L64").Value = 0
$ws_WVR.Range("N64").ClearContents()

$ws_WVR.Range("H67").Value = 0
$ws_WVR.Range("J67").Value = 0
$ws_WVR.Range("L67").Value = 0
$ws_WVR.Range("N67").ClearContents()

$ws_WVR.Range("H113").Value = 1061.5186
$ws_WVR.Range("I113").Value = 1040.0555
$ws_WVR.Range("K113").Value = 3120.1665
$ws_WVR.Range("M113").Value = -950.1664999999998

$ws_WVR.Range("H122").Value = 43581.188
$ws_WVR.Range("I122").Value = 49625.895
$ws_WVR.Range("J122").Value = 1268.25
$ws_WVR.Range("K122").Value = 148877.685
$ws_WVR.Range("L122").Value = 3804.75
$ws_WVR.Range("M122").Value = -146427.685
$ws_WVR.Range("N122").Value = -8704.75

$ws_WVR.Range("H126").Value = 5841
$ws_WVR.Range("I126").Value = 4600
$ws_WVR.Range("J126").Value = 6668.3335
$ws_WVR.Range("K126").Value = 13800
$ws_WVR.Range("L126").Value = 20005.0005
$ws_WVR.Range("M126").Value = -11330
$ws_WVR.Range("N126").Value = -24945.0005

$ws_WVR.Range("H136").Value = 29821.625
$ws_WVR.Range("I136").Value = 35096
$ws_WVR.Range("K136").Value = 105288
$ws_WVR.Range("M136").Value = -102738
